# Hide slides 3 through 18 (inclusive) in the deck.
# This mirrors PowerPoint's "Hide Slide" command, which on save is persisted
# as show="0" on the <p:sld> element and is surfaced through the object
# model as Slide.SlideShowTransition.Hidden (msoTrue/msoFalse).
$p = $ppt.ActivePresentation

for ($i = 3; $i -le 18; $i++) {
    $slide = $p.Slides.Item($i)
    $slide.SlideShowTransition.Hidden = $true
}
